$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the very start of the document
# (wrapping the opening of paragraph 1). It needs to move to just before
# "ion." near the end of the final paragraph's last run ("...information.").
#
# Locate "information." and split the insertion point 8 characters in,
# i.e. right after "informat" / right before "ion.".
$rng = $d.Content
$found = $rng.Find.Execute("information.", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

$splitPoint = $rng.Start + 8
$bmRange = $d.Range($splitPoint, $splitPoint)

# Re-adding a bookmark with the same name relocates it (removing the old
# occurrence) rather than creating a duplicate, which is exactly what we
# want here.
$d.Bookmarks.Add("_GoBack", $bmRange)
